# R-任务-任务配置 (task config) resource update, per "update resource by parse tool".
# Reworks the sample task rows 5-7 (difficulty-tiered "complete level" tasks)
# and removes the stray task-id values left in A8:A10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: easy-tier task (was "初级"/"2次") -> "简单" / "1次"
$ws.Range("G5").Value = "完成关卡(简单)"
$ws.Range("H5").Value = "完成关卡1次"

# Row 6: normal-tier task, renumber the linked level id, keep display text,
# but completion count moves to "5次"
$ws.Range("E6").Value = 10102
$ws.Range("G6").Value = "完成关卡(普通)"
$ws.Range("H6").Value = "完成关卡5次"

# Row 7: new hard-tier task row - pick up G6:H6's cell formatting first,
# then overwrite the values
$ws.Range("G6:H6").Copy()
$ws.Range("G7:H7").PasteSpecial(-4122)
$ws.Range("E7").Value = 10103
$ws.Range("G7").Value = "完成关卡(困难)"
$ws.Range("H7").Value = "完成关卡10次"

# Rows 8-10 lost their stray leading task-id values (B/C stay as-is)
$ws.Range("A8").ClearContents()
$ws.Range("A9").ClearContents()
$ws.Range("A10").ClearContents()

# Update the saved selection to match where the editor left off
$ws.Range("G8").Select()
